$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03400945102512423
$ws.Range("C2").Value = 1.399455635964034
$ws.Range("D2").Value = 6.170139677499137
$ws.Range("E2").Value = 2.483976585537621
$ws.Range("F2").Value = 2.542192780186525
$ws.Range("G2").Value = 22
